$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (old B -> C, old C -> D)
$ws.Columns.Item(2).Insert()

# The insert operation copies formatting from the column to the left (A) onto
# the new column B data rows; the target layout has no style on those cells,
# so clear it.
$ws.Range("B2:B20").ClearFormats()

# New header for column B: "segments", styled like the other header cells
# (copy format from C1, which holds the old "PercActivations" header).
$ws.Cells.Item(1, 3).Copy() | Out-Null
$ws.Cells.Item(1, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 2).Value = "segments"

# Rows 2..20: move the segment-name text from column A into the new column B,
# then replace column A with the numeric (0-based) segment index, keeping A's
# existing style.
for ($row = 2; $row -le 20; $row++) {
    $segmentName = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 2).Value = $segmentName
    $ws.Cells.Item($row, 1).Value = $row - 2
}

$excel.CutCopyMode = 0
